$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D/E columns are treated as text so values such as "1.006" or
# "28.678.61" are not auto-converted into numbers by Excel.
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.678.61'
$ws.Range("E2").Value = '  +2.29%  '

$ws.Range("D3").Value = '1.872.57'
$ws.Range("E3").Value = '  +2.25%  '

$ws.Range("D4").Value = '1.006'
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").Value = '326.14'
$ws.Range("E5").Value = '  -0.57%  '

$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.15%  '

$ws.Range("D7").Value = '0.4642'
$ws.Range("E7").Value = '  +0.45%  '

$ws.Range("D8").Value = '0.3902'
$ws.Range("E8").Value = '  +1.29%  '

$ws.Range("D9").Value = '0.07909'
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("D10").Value = '0.9744'
$ws.Range("E10").Value = '  +1.59%  '

$ws.Range("D11").Value = '22.12'
$ws.Range("E11").Value = '  +1.22%  '

$ws.Range("D12").Value = '1.832.84'
$ws.Range("E12").Value = '  -4.79%  '

$ws.Range("D13").Value = '5.718'
$ws.Range("E13").Value = '  +1.25%  '

$ws.Range("D14").Value = '6.976'
$ws.Range("E14").Value = '  +1.45%  '

$ws.Range("D15").Value = '0.06983'
$ws.Range("E15").Value = '  +2.57%  '

$ws.Range("D16").Value = '88.32'
$ws.Range("E16").Value = '  +1.87%  '

$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.33%  '

$ws.Range("D18").Value = '0.00001007'
$ws.Range("E18").Value = '  +1.06%  '

$ws.Range("D19").Value = '16.84'
$ws.Range("E19").Value = '  +1.50%  '

$ws.Range("D20").Value = '1.005'
$ws.Range("E20").Value = '  +0.23%  '

$ws.Range("D21").Value = '28.675.80'
$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("D22").Value = '5.316'
$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("B24").Value = 'Toncoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D24").Value = '2.118'
$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("B25").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C25").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D25").Value = '2.107.00'
$ws.Range("E25").Value = '  +0.14%  '

$ws.Range("D26").Value = '152.87'
$ws.Range("E26").Value = '  -0.29%  '

$ws.Range("D27").Value = '19.28'
$ws.Range("E27").Value = '  +0.43%  '

$ws.Range("D28").Value = '5.734'
$ws.Range("E28").Value = '  +0.90%  '

$ws.Range("D29").Value = '1.989'
$ws.Range("E29").Value = '  +0.84%  '

$ws.Range("E30").Value = '  +2.09%  '

$ws.Range("D31").Value = '0.09350'
$ws.Range("E31").Value = '  +1.00%  '

$ws.Range("D32").Value = '0.9264'
$ws.Range("E32").Value = '  -0.71%  '

$ws.Range("D33").Value = '5.281'
$ws.Range("E33").Value = '  -0.22%  '

$ws.Range("D34").Value = '1.342'
$ws.Range("E34").Value = '  +1.91%  '

$ws.Range("D35").Value = '3.346'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("D36").Value = '0.05819'
$ws.Range("E36").Value = '  -0.78%  '

$ws.Range("D37").Value = '0.02098'
$ws.Range("E37").Value = '  -1.98%  '

$ws.Range("D38").Value = '1.144'

$ws.Range("D39").Value = '7.772'
$ws.Range("E39").Value = '  +1.11%  '

$ws.Range("D40").Value = '0.5663'
$ws.Range("E40").Value = '  +1.45%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.1786'
$ws.Range("E41").Value = '  +1.33%  '

$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").Value = '9.901'
$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("D43").Value = '0.07246'
$ws.Range("E43").Value = '  +3.29%  '

$ws.Range("D44").Value = '11.71'
$ws.Range("E44").Value = '  +1.19%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.5320'
$ws.Range("E45").Value = '  +1.06%  '

$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '1.164'
$ws.Range("E46").Value = '  -3.60%  '

$ws.Range("D47").Value = '1.831'
$ws.Range("E47").Value = '  +0.11%  '

$ws.Range("D48").Value = '113.29'
$ws.Range("E48").Value = '  +1.10%  '

$ws.Range("D49").Value = '2.036'
$ws.Range("E49").Value = '  -5.18%  '

$ws.Range("D50").Value = '2.369'
$ws.Range("E50").Value = '  +1.75%  '

$ws.Range("E51").Value = '  +0.31%  '

# Restore the default (Normal) style on the value range so that the
# explicit text number format introduced above does not linger as a
# visible style change on the cells.
$valueRange.Style = "Normal"
